$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.632.43"
$ws.Range("E2").Value = "  -3.27%  "

$ws.Range("D3").Value = "2.085.96"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "'344.80"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").Value = "'0.5159"
$ws.Range("E7").Value = "  -2.15%  "

$ws.Range("D8").Value = "'0.4386"
$ws.Range("E8").Value = "  -2.98%  "

$ws.Range("D9").Value = "'0.09182"
$ws.Range("E9").Value = "  +1.44%  "

$ws.Range("D10").Value = "'51.70"
$ws.Range("E10").Value = "  -4.01%  "

$ws.Range("D11").Value = "'1.172"
$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("D12").Value = "'25.38"
$ws.Range("E12").Value = "  +3.68%  "

$ws.Range("D13").Value = "2.087.36"
$ws.Range("E13").Value = "  -1.27%  "

# Row 14/15: Chainlink and Polkadot swap positions
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.723"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'8.176"
$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").Value = "'0.00001158"
$ws.Range("E17").Value = "  -1.99%  "

$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("D19").Value = "'21.11"
$ws.Range("E19").Value = "  +9.03%  "

$ws.Range("D20").Value = "'0.06640"
$ws.Range("E20").Value = "  -1.35%  "

$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("D22").Value = "'6.200"
$ws.Range("E22").Value = "  -2.10%  "

$ws.Range("D23").Value = "29.730.88"
$ws.Range("E23").Value = "  -3.15%  "

$ws.Range("D24").Value = "'12.67"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = "  -3.85%  "

$ws.Range("D26").Value = "2.333.52"
$ws.Range("E26").Value = "  -1.18%  "

$ws.Range("D27").Value = "'21.88"
$ws.Range("E27").Value = "  -2.40%  "

$ws.Range("D28").Value = "'163.00"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("D29").Value = "'2.516"
$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("D30").Value = "'132.77"
$ws.Range("E30").Value = "  -2.46%  "

$ws.Range("D31").Value = "'1.145"
$ws.Range("E31").Value = "  -4.15%  "

$ws.Range("D32").Value = "'0.1052"
$ws.Range("E32").Value = "  -2.12%  "

$ws.Range("D33").Value = "'1.626"
$ws.Range("E33").Value = "  -0.64%  "

$ws.Range("D34").Value = "'6.190"
$ws.Range("E34").Value = "  -2.73%  "

$ws.Range("D35").Value = "'3.938"
$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("D36").Value = "'6.102"
$ws.Range("E36").Value = "  +2.91%  "

$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("E38").Value = "  -2.59%  "

$ws.Range("D39").Value = "'0.06720"
$ws.Range("E39").Value = "  -1.84%  "

$ws.Range("D40").Value = "'0.2272"
$ws.Range("E40").Value = "  -1.92%  "

$ws.Range("D41").Value = "'12.45"
$ws.Range("E41").Value = "  -1.40%  "

$ws.Range("D42").Value = "'0.6844"
$ws.Range("E42").Value = "  -0.57%  "

$ws.Range("D43").Value = "'1.288"
$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").Value = "'0.6639"
$ws.Range("E44").Value = "  +3.33%  "

$ws.Range("D45").Value = "'14.15"
$ws.Range("E45").Value = "  -4.57%  "

$ws.Range("D46").Value = "'2.295"
$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("D47").Value = "'3.621"
$ws.Range("E47").Value = "  -3.70%  "

$ws.Range("D48").Value = "'1.216"
$ws.Range("E48").Value = "  -3.28%  "

# Row 49/50: BabyDogeCoin and Aave swap positions
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'81.51"
$ws.Range("E49").Value = "  -1.70%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000331"
$ws.Range("E50").Value = "  -7.75%  "

$ws.Range("D51").Value = "'1.166"
$ws.Range("E51").Value = "  -2.56%  "
